# Reorder the score table: move the "HIK" row (originally the last data
# row, row 7) up to become the first data row (row 2), shifting the
# other rows (originally 2-6) down by one (to 3-7).
#
# Columns: A = mål (label), B = score, C = 95%ki, D = percentil

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the last row (HIK) before it gets overwritten by the shift below.
$lastRowValues = @(
    $ws.Cells.Item(7, 1).Value2,
    $ws.Cells.Item(7, 2).Value2,
    $ws.Cells.Item(7, 3).Value2,
    $ws.Cells.Item(7, 4).Value2
)

# Shift rows 2..6 down into rows 3..7 (walk bottom-up so we never
# clobber a row before it has been read).
for ($r = 6; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r + 1, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# Write the stashed HIK row into the new first data row (row 2).
$ws.Cells.Item(2, 1).Value = $lastRowValues[0]
$ws.Cells.Item(2, 2).Value = $lastRowValues[1]
$ws.Cells.Item(2, 3).Value = $lastRowValues[2]
$ws.Cells.Item(2, 4).Value = $lastRowValues[3]

# Update the active selection to match the author's final cursor position.
$ws.Range("H9").Select()
